$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 322, shifting existing rows 322-354 down to 323-355
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new data record
$ws.Cells.Item(322, 1).Value = 4
$ws.Cells.Item(322, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(322, 3).Value = "Los Lagos"
$ws.Cells.Item(322, 4).Value = 44918
$ws.Cells.Item(322, 5).Value = 10
$ws.Cells.Item(322, 6).Value = 100112021
$ws.Cells.Item(322, 7).Value = "Ají"
$ws.Cells.Item(322, 8).Value = "Inferno"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 90
$ws.Cells.Item(322, 11).Value = 22000
$ws.Cells.Item(322, 12).Value = 22000
$ws.Cells.Item(322, 13).Value = 22000
$ws.Cells.Item(322, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(322, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(322, 16).Value = 2200
$ws.Cells.Item(322, 17).Value = 10
$ws.Cells.Item(322, 18).Value = "Hortaliza"
